$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 211.66667
$ws.Range("I28").Value = 211.66667
$ws.Range("J28").Value = 0
$ws.Range("K28").Value = 211.66667
$ws.Range("L28").Value = 0
$ws.Range("M28").Value = 273.33333
$ws.Range("N28").ClearContents()
$ws.Range("H53").Value = 347.625
$ws.Range("I53").Value = 272.25
$ws.Range("K53").Value = 272.25
$ws.Range("M53").Value = 364.75
$ws.Range("H58").Value = 2396.6
$ws.Range("I58").Value = 1999
$ws.Range("J58").Value = 2496
$ws.Range("K58").Value = 5997
$ws.Range("L58").Value = 7488
$ws.Range("M58").Value = -5847
$ws.Range("N58").Value = -7788
$ws.Range("H100").Value = 5500
$ws.Range("I100").Value = 10000
$ws.Range("J100").Value = 1000
$ws.Range("K100").Value = 10000
$ws.Range("L100").Value = 1000
$ws.Range("M100").Value = -9459
$ws.Range("N100").Value = -2082
$ws.Range("H101").Value = 2098.6
$ws.Range("J101").Value = 498.66666
$ws.Range("L101").Value = 1495.99998
$ws.Range("N101").Value = -4739.999980000001
$ws.Range("H137").Value = 1666.3334
$ws.Range("J137").Value = 1624.5
$ws.Range("L137").Value = 4873.5
$ws.Range("N137").Value = -9973.5

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 14692
$ws.Range("I32").Value = 14692
$ws.Range("K32").Value = 14692
$ws.Range("M32").Value = -14405
$ws.Range("H45").Value = 5499.6665
$ws.Range("I45").Value = 5499.6665
$ws.Range("K45").Value = 5499.6665
$ws.Range("M45").Value = -5122.6665
$ws.Range("H61").Value = 1819.4
$ws.Range("I61").Value = 2024.5
$ws.Range("J61").Value = 999
$ws.Range("K61").Value = 2024.5
$ws.Range("L61").Value = 999
$ws.Range("M61").Value = -1812.5
$ws.Range("N61").Value = -1423
$ws.Range("H74").Value = 16691.9
$ws.Range("I74").Value = 16363.625
$ws.Range("J74").Value = 18005
$ws.Range("K74").Value = 16363.625
$ws.Range("L74").Value = 18005
$ws.Range("M74").Value = -15489.625
$ws.Range("N74").Value = -19753
$ws.Range("H77").Value = 16691.9
$ws.Range("I77").Value = 16363.625
$ws.Range("J77").Value = 18005
$ws.Range("K77").Value = 81818.125
$ws.Range("L77").Value = 90025
$ws.Range("M77").Value = -77450.125
$ws.Range("N77").Value = -98761
$ws.Range("H132").Value = 2894.3076
$ws.Range("I132").Value = 2718.9167
$ws.Range("K132").Value = 8156.750100000001
$ws.Range("M132").Value = -5626.750100000001
$ws.Range("H136").Value = 1819.4
$ws.Range("I136").Value = 2024.5
$ws.Range("J136").Value = 999
$ws.Range("K136").Value = 6073.5
$ws.Range("L136").Value = 2997
$ws.Range("M136").Value = -3523.5
$ws.Range("N136").Value = -8097

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H61").Value = 0
$ws.Range("J61").Value = 0
$ws.Range("L61").Value = 0
$ws.Range("N61").ClearContents()

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2384.5833
$ws.Range("I31").Value = 2282
$ws.Range("J31").Value = 2528.2
$ws.Range("K31").Value = 2282
$ws.Range("L31").Value = 2528.2
$ws.Range("M31").Value = -1987
$ws.Range("N31").Value = -3118.2
$ws.Range("H34").Value = 2384.5833
$ws.Range("I34").Value = 2282
$ws.Range("J34").Value = 2528.2
$ws.Range("K34").Value = 2282
$ws.Range("L34").Value = 2528.2
$ws.Range("M34").Value = -2080
$ws.Range("N34").Value = -2932.2
$ws.Range("H39").Value = 3050.5
$ws.Range("I39").Value = 3050.5
$ws.Range("J39").Value = 0
$ws.Range("K39").Value = 3050.5
$ws.Range("L39").Value = 0
$ws.Range("M39").Value = -2659.5
$ws.Range("N39").ClearContents()
$ws.Range("H49").Value = 3050.5
$ws.Range("I49").Value = 3050.5
$ws.Range("J49").Value = 0
$ws.Range("K49").Value = 3050.5
$ws.Range("L49").Value = 0
$ws.Range("M49").Value = -2868.5
$ws.Range("N49").ClearContents()
$ws.Range("H69").Value = 0
$ws.Range("I69").Value = 0
$ws.Range("K69").Value = 0
$ws.Range("M69").ClearContents()
$ws.Range("H72").Value = 0
$ws.Range("I72").Value = 0
$ws.Range("K72").Value = 0
$ws.Range("M72").ClearContents()
$ws.Range("H105").Value = 5449
$ws.Range("J105").Value = 4748.3335
$ws.Range("L105").Value = 4748.3335
$ws.Range("N105").Value = -8242.333500000001
$ws.Range("H141").Value = 56499.5
$ws.Range("J141").Value = 56499.5
$ws.Range("L141").Value = 56499.5
$ws.Range("N141").Value = -66859.5

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H23").Value = 1999.6666
$ws.Range("I23").Value = 1999.6666
$ws.Range("K23").Value = 5998.9998
$ws.Range("M23").Value = -5763.9998
$ws.Range("H25").Value = 1500
$ws.Range("J25").Value = 0
$ws.Range("L25").Value = 0
$ws.Range("N25").ClearContents()
$ws.Range("H30").Value = 1500
$ws.Range("J30").Value = 0
$ws.Range("L30").Value = 0
$ws.Range("N30").ClearContents()
$ws.Range("H37").Value = 149800
$ws.Range("J37").Value = 149800
$ws.Range("L37").Value = 449400
$ws.Range("N37").Value = -449624
$ws.Range("H97").Value = 629
$ws.Range("I97").Value = 604.13336
$ws.Range("K97").Value = 1812.40008
$ws.Range("M97").Value = -1316.40008
$ws.Range("H98").Value = 2470.625
$ws.Range("I98").Value = 3368
$ws.Range("J98").Value = 1573.25
$ws.Range("K98").Value = 10104
$ws.Range("L98").Value = 4719.75
$ws.Range("M98").Value = -8606
$ws.Range("N98").Value = -7715.75
$ws.Range("H131").Value = 2504.25
$ws.Range("I131").Value = 1512.5
$ws.Range("K131").Value = 4537.5
$ws.Range("M131").Value = 502.5
$ws.Range("H136").Value = 2000
$ws.Range("I136").Value = 2000
$ws.Range("K136").Value = 6000
$ws.Range("M136").Value = -900
$ws.Range("H139").Value = 2516
$ws.Range("I139").Value = 999
$ws.Range("J139").Value = 4033
$ws.Range("K139").Value = 2997
$ws.Range("L139").Value = 12099
$ws.Range("M139").Value = 2143
$ws.Range("N139").Value = -22379
$ws.Range("H140").Value = 3206
$ws.Range("I140").Value = 3206
$ws.Range("K140").Value = 9618
$ws.Range("M140").Value = -4438

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H99").Value = 7710.857
$ws.Range("I99").Value = 2662.6667
$ws.Range("K99").Value = 2662.6667
$ws.Range("M99").Value = -416.6667000000002
$ws.Range("H138").Value = 20000
$ws.Range("J138").Value = 20000
$ws.Range("L138").Value = 20000
$ws.Range("N138").Value = -30280

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 5674.4165
$ws.Range("I46").Value = 1644.8182
$ws.Range("K46").Value = 1644.8182
$ws.Range("M46").Value = -1456.8182
$ws.Range("I82").Value = 3784.2
$ws.Range("J82").Value = 3333
$ws.Range("K82").Value = 3784.2
$ws.Range("L82").Value = 3333
$ws.Range("M82").Value = -3423.2
$ws.Range("N82").Value = -4055
$ws.Range("I85").Value = 3784.2
$ws.Range("J85").Value = 3333
$ws.Range("K85").Value = 3784.2
$ws.Range("L85").Value = 3333
$ws.Range("M85").Value = -2536.2
$ws.Range("N85").Value = -5829
$ws.Range("H93").Value = 3278.8333
$ws.Range("I93").Value = 3278.8333
$ws.Range("K93").Value = 3278.8333
$ws.Range("M93").Value = -2030.8333

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 6618.222
$ws.Range("I136").Value = 4652.143
$ws.Range("K136").Value = 13956.429
$ws.Range("M136").Value = -11406.429
